# Atualização de bases das ligas, do dia: 30-05-2024 às 12:21
#
# Several fixtures had their two rows written in the wrong order
# relative to one another. This swaps the match data (columns B
# through AD) between the two rows of each affected pair, while
# leaving column A (the sequential row index) untouched. Cells whose
# value is identical between the two rows are left untouched so their
# original on-disk representation is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(18, 19),
    @(54, 55),
    @(58, 59),
    @(63, 64),
    @(67, 68)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        if ($val1 -ceq $val2) {
            continue
        }

        $cell1.Formula = $val2
        $cell2.Formula = $val1
    }
}
